$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the stale C:E columns (Count / First Seen Time / Last Seen Time)
# entirely, keeping A:B formatting intact.
$ws.Range("C1:E3").Clear()

# Header row (A1/B1 keep their existing bold/border header style)
$ws.Range("A1").Value = "Name"
$ws.Range("B1").Value = "Image_Path"

# Data rows - clear old contents first (row 2/3 had Date/Count/etc in B:E)
$ws.Range("A2:B7").ClearContents()

$ws.Range("A2").Value = "Malindha"
$ws.Range("B2").Value = "person2.jpg"

$ws.Range("A3").Value = "Lalithya"
$ws.Range("B3").Value = "person3.jpg"

$ws.Range("A4").Value = "Anuradha"
$ws.Range("B4").Value = "person-4.jpg"

$ws.Range("A5").Value = "Malinga"
$ws.Range("B5").Value = "person6.jpg"

$ws.Range("A6").Value = "Hasith"
$ws.Range("B6").Value = "person8.jpg"

# "0012" must stay text (leading zeros), so force text format before
# assigning the value to avoid Excel coercing it to the number 12, then
# drop the number-format override again so the cell keeps the default
# (unstyled) look of the other data rows.
$ws.Range("A7").NumberFormat = "@"
$ws.Range("A7").Value = "0012"
$ws.Range("A7").ClearFormats()
$ws.Range("B7").Value = "0012.jpg"
